$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection from D8 to E10
$ws.Range("E10").Select()
